# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Also recalc std/mean and write s_vals (handled automatically by Excel
# recalculation once the underlying K values are corrected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 3
    4  = 4
    5  = 3
    6  = 5
    7  = 5
    8  = 5
    9  = 1
    10 = 3
    11 = 0
    12 = 4
    13 = 2
    14 = 6
    15 = 3
    16 = 6
    17 = 3
    18 = 1
    19 = 2
    20 = 2
    21 = 5
    22 = 6
    23 = 2
    24 = 5
    25 = 3
    26 = 4
    27 = 9
    28 = 2
    29 = 7
    30 = 6
    31 = 4
    32 = 5
    33 = 4
    34 = 2
    35 = 3
    36 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}

$excel.CalculateFull()
$wb.Save()
